$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: fill in previously-empty BOM line (ceramic cap) ---
# Write Manufacturer / Part number / Part description first so that the
# shared-string table gets populated in the same order the source file uses.
$ws.Range("D12").Value = "Multicomp"
$ws.Range("F12").Value = "MC0805B224K500CT"
$ws.Range("B12").Value = "0,22uF 0805 50V X7R"
$ws.Range("H12").Value = 2320842
$ws.Range("L12").Value = 0.0843

# --- Row 13: fill in previously-empty BOM line (resistor) ---
$ws.Range("B13").Value = "120R 0603"
$ws.Range("D13").Value = "Vishay"
$ws.Range("F13").Value = "CRCW0603120RFKEA "
$ws.Range("H13").Value = 1652832

# --- Row 25: new BOM line (PTC fuse) ---
$ws.Range("B25").Value = "PTC "
$ws.Range("D25").Value = "Littlefuse"
$ws.Range("F25").Value = "1.1A/1.95A 1812L110/33MR"
$ws.Range("H25").Value = 1822213
$ws.Range("L25").Value = 0.703

# --- Row 26: new BOM line (resistor) ---
$ws.Range("B26").Value = "215k 0603"
$ws.Range("D26").Value = "Vishay"
$ws.Range("F26").Value = "CRCW0603215KFKEA"
$ws.Range("H26").Value = 2138528
$ws.Range("L26").Value = 0.0174

# --- Row 27: new BOM line (resistor) ---
$ws.Range("B27").Value = "5k6 0603"
$ws.Range("D27").Value = "Multicomp"
$ws.Range("F27").Value = "mcmr06x5601ftl"
$ws.Range("H27").Value = 2073537
$ws.Range("L27").Value = 0.0087

# --- View state: zoom to 85% and move the active selection to H5 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("H5").Select()
